$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.246.04'
$ws.Range('E2').Value = '  -1.89%  '
$ws.Range('D3').Value = '2.500.23'
$ws.Range('E3').Value = '  -4.87%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.60'
$ws.Range('E5').Value = '  -2.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.65'
$ws.Range('E6').Value = '  +1.35%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -2.74%  '
$ws.Range('D9').Value = '2.497.57'
$ws.Range('E9').Value = '  -5.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.138'
$ws.Range('E10').Value = '  -1.51%  '
$ws.Range('E11').Value = '  -0.36%  '
$ws.Range('E12').Value = '  -4.52%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.08'
$ws.Range('E13').Value = '  -2.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.40'
$ws.Range('E14').Value = '  -4.55%  '
$ws.Range('D15').Value = '2.939.52'
$ws.Range('E15').Value = '  -5.46%  '
$ws.Range('E16').Value = '  -4.51%  '
$ws.Range('D17').Value = '66.144.63'
$ws.Range('E17').Value = '  -1.55%  '
$ws.Range('D18').Value = '2.477.26'
$ws.Range('E18').Value = '  -6.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.19'
$ws.Range('E19').Value = '  -6.88%  '
$ws.Range('E20').Value = '  -5.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '344.91'
$ws.Range('E21').Value = '  -3.36%  '
$ws.Range('E22').Value = '  -3.33%  '
$ws.Range('E23').Value = '  -2.52%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.92'
$ws.Range('E25').Value = '  -0.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '69.10'
$ws.Range('E26').Value = '  -0.83%  '
$ws.Range('E27').Value = '  -4.71%  '
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').Value = '2.623.97'
$ws.Range('E29').Value = '  -5.18%  '
$ws.Range('E30').Value = '  -4.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '524.94'
$ws.Range('E31').Value = '  -4.06%  '
$ws.Range('E32').Value = '  +1.74%  '
$ws.Range('E33').Value = '  -2.86%  '
$ws.Range('E34').Value = '  -3.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.130'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '157.60'
$ws.Range('E37').Value = '  +0.65%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.44'
$ws.Range('E38').Value = '  -4.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.48'
$ws.Range('E39').Value = '  -2.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.34'
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('E41').Value = '  -3.84%  '
$ws.Range('E42').Value = '  -2.69%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.04'
$ws.Range('E43').Value = '  -3.63%  '
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('E45').Value = '  +0.32%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '146.40'
$ws.Range('E46').Value = '  -4.23%  '
$ws.Range('E47').Value = '  -4.51%  '
$ws.Range('E48').Value = '  -3.69%  '
$ws.Range('E49').Value = '  +1.20%  '
$ws.Range('D50').Value = '0.0₆0268'
$ws.Range('E50').Value = '  -10.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0750'
$ws.Range('E51').Value = '  -2.69%  '
